$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B32: change from text "2" to numeric 2
$ws.Range("B32").Value = 2

# Add new row 33
$ws.Range("A33").Value = "Sunsi Wu"
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "4"
$ws.Range("C33").Value = "will"
$ws.Range("D33").Value = "ACK"
$ws.Range("E33").Value = "OTH"
$ws.Range("F33").Value = "d4ad31e6-de82-4ee8-af90-c18d97ed2c36"
$ws.Range("G33").Value = "Bk7wvW-C-_annotated.xlsx"
$ws.Range("H33").Value = "We will update our paper very soon."
